$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the two completed tasks from "Doing:" column (E) to the end of the
# "Done:" column (G) - i.e. mark them as finished on the scrum board.
$taskReunir = $ws.Range("E4").Value()
$taskDecompor = $ws.Range("E5").Value()

# Clear the old locations in column E
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()

# Append them at the bottom of the Done column
$ws.Range("G7").Value = $taskReunir
$ws.Range("G8").Value = $taskDecompor

# Update the view: scroll so column C is leftmost and select G13, matching
# the cursor position recorded when the workbook was last saved.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G13").Select()
